# BOT; UPDATE DATA
# Adds one more day (2020-04-18, serial 43939) of figures to the three
# data sheets ("all", "kobe", "other"), shifting the trailing
# label/footer rows down by one row, and corrects one previously
# reported figure (kobe!D65).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "all": insert new data row 11 (old rows 11-12 shift to 12-13)
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Rows.Item(11).Insert()

# Match formatting of the previous day's row for the new row, then
# overwrite the columns whose formatting differs from a straight
# copy-down (D:H keep the "normal" style used by earlier rows).
$wsAll.Range("D9:H9").Copy() | Out-Null
$wsAll.Range("D11:H11").PasteSpecial(-4122) | Out-Null

$wsAll.Range("A11").Value = 43939
$wsAll.Range("B11").Value = 179
$wsAll.Range("C11").Value = 146
$wsAll.Range("D11").Value = 108
$wsAll.Range("E11").Value = 100
$wsAll.Range("F11").Value = 8
$wsAll.Range("G11").Value = 2
$wsAll.Range("H11").Value = 36

$wsAll.Range("B11:H11").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "kobe": correct D65, then insert new data row 66
# (old row 66 shifts to 67)
# ---------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D65").Value = 12

$wsKobe.Rows.Item(66).Insert()

# F:J of the new row use the "normal" style (style index 2) which is
# not otherwise present on this sheet yet - copy it in from sheet
# "all" where it is already used, then fix C66 back to style 7.
$wsAll.Range("D9:H9").Copy() | Out-Null
$wsKobe.Range("F66:J66").PasteSpecial(-4122) | Out-Null
$wsKobe.Range("C64").Copy() | Out-Null
$wsKobe.Range("C66").PasteSpecial(-4122) | Out-Null

$wsKobe.Range("A66").Value = 43939
$wsKobe.Range("B66").Value = 16
$wsKobe.Range("C66").Formula = "=C65+B66"
$wsKobe.Range("D66").Value = 6
$wsKobe.Range("E66").Value = 179
$wsKobe.Range("F66").Value = 101
$wsKobe.Range("G66").Value = 94
$wsKobe.Range("H66").Value = 7
$wsKobe.Range("I66").Value = 2
$wsKobe.Range("J66").Value = 34

$wsKobe.Range("F15").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "other": insert new data row 41
# (old rows 41-42 shift to 42-43)
# ---------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Rows.Item(41).Insert()

$wsOther.Range("A40:I40").Copy() | Out-Null
$wsOther.Range("A41:I41").PasteSpecial(-4122) | Out-Null

$wsOther.Range("A41").Value = 43939
$wsOther.Range("B41").Value = 0
$wsOther.Range("C41").Value = 9
$wsOther.Range("D41").Value = 7
$wsOther.Range("E41").Value = 6
$wsOther.Range("F41").Value = 1
$wsOther.Range("G41").Value = 0
$wsOther.Range("H41").Value = 2

$wsOther.Range("D48").Select() | Out-Null

# ---------------------------------------------------------------
# Keep "all" as the active/selected sheet, matching the workbook's
# original tab selection.
# ---------------------------------------------------------------
$wsAll.Activate() | Out-Null
